# Updates cryptocurrency price (D) and 1h volume-change (E) columns
# to the latest scraped snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.300.53"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "1.595.47"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'211.55"
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  -0.18%  "
$ws.Range("D10").Value = "'19.01"
$ws.Range("E10").Value = "  +0.01%  "
$ws.Range("E11").Value = "  +0.96%  "
$ws.Range("D12").Value = "1.820.84"
$ws.Range("E12").Value = "  +0.43%  "
$ws.Range("D13").Value = "1.600.23"
$ws.Range("E13").Value = "  +0.66%  "
$ws.Range("E14").Value = "  -0.78%  "
$ws.Range("E15").Value = "  -1.47%  "
$ws.Range("D16").Value = "'63.39"
$ws.Range("E16").Value = "  -0.37%  "
$ws.Range("D17").Value = "26.288.39"
$ws.Range("E17").Value = "  +0.46%  "
$ws.Range("D18").Value = "'229.48"
$ws.Range("E18").Value = "  +6.75%  "
$ws.Range("E19").Value = "  +3.88%  "
$ws.Range("E20").Value = "  -0.45%  "
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("D23").Value = "'2.17"
$ws.Range("E23").Value = "  +2.40%  "
$ws.Range("E24").Value = "  -1.08%  "
$ws.Range("D25").Value = "'146.37"
$ws.Range("E25").Value = "  +1.25%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("E28").Value = "  -0.17%  "
$ws.Range("D29").Value = "'15.39"
$ws.Range("E29").Value = "  +2.09%  "
$ws.Range("D30").Value = "'0.0494"
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("E31").Value = "  -0.13%  "
$ws.Range("D32").Value = "1.489.31"
$ws.Range("E32").Value = "  +5.04%  "
$ws.Range("E33").Value = "  +0.88%  "
$ws.Range("E34").Value = "  -0.94%  "
$ws.Range("E35").Value = "  -0.42%  "
$ws.Range("E36").Value = "  +0.25%  "
$ws.Range("D37").Value = "'0.568"
$ws.Range("E37").Value = "  -2.97%  "
$ws.Range("E38").Value = "  -0.47%  "
$ws.Range("E39").Value = "  -0.99%  "
$ws.Range("E40").Value = "  -1.52%  "
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("E42").Value = "  +1.25%  "
$ws.Range("D43").Value = "'0.927"
$ws.Range("E43").Value = "  -1.98%  "
$ws.Range("D44").Value = "1.733.41"
$ws.Range("E44").Value = "  +0.52%  "
$ws.Range("E45").Value = "  -0.82%  "
$ws.Range("D46").Value = "'60.30"
$ws.Range("E46").Value = "  -1.29%  "
$ws.Range("D47").Value = "'88.44"
$ws.Range("E47").Value = "  +1.23%  "
$ws.Range("E48").Value = "  -0.63%  "
$ws.Range("E49").Value = "  -0.21%  "
$ws.Range("E50").Value = "  -0.63%  "
$ws.Range("E51").Value = "  +0.08%  "

# The price strings above (e.g. "211.55") parse as plain numbers, so the
# leading apostrophe forces Excel to store them as text like the source
# sheet does; reset the style so no stray quote-prefix formatting sticks.
$ws.Range("D5").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
